# Add participant 240M_FM to the Alpha matlab genotypes worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data occupies A1:B81 (header row + 80 participants).
# Apply a left-aligned style to the genotype column (B) for all existing rows -
# this mirrors the new cellXfs entry (left horizontal alignment) that gets
# referenced by every cell in column B.
$ws.Range("B1:B81").HorizontalAlignment = -4131

# Record the column-level formatting for column B as well so new cells in
# that column pick up the same default width/style.
$ws.Columns("B").ColumnWidth = 9.140625

# Append the new participant as row 82.
$ws.Cells.Item(82, 1).Value = "240M_FM"
$ws.Cells.Item(82, 2).Value = 0
$ws.Range("B82").HorizontalAlignment = -4131

# Restore the author's view state (the sheet was scrolled/selected around
# the newly added row when the workbook was saved).
$ws.Range("I77").Select() | Out-Null

Write-Output "Added participant 240M_FM as row 82"
